$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 15
$ws.Cells.Item($row - 1, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = 42620.886458333334
$ws.Cells.Item($row, 2).Value = -10
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 41
$ws.Cells.Item($row, 5).Value = 55
$ws.Cells.Item($row, 6).Value = 70
$ws.Cells.Item($row, 7).Value = 10803
$ws.Cells.Item($row, 8).Value = 17151
$ws.Cells.Item($row, 9).Value = 1835
$ws.Cells.Item($row, 10).Value = 220
$ws.Cells.Item($row, 11).Value = 166
$ws.Cells.Item($row, 12).Value = 5
$ws.Cells.Item($row, 13).Value = 12
$ws.Cells.Item($row, 14).Value = "Bag"
